# cryptos.xlsx refresh -- "Updated symbol list ... with GitHub Actions"
#
# Every touched cell on the sheet is stored as text (coin names, URLs,
# decimal prices and percentage strings are all plain inline strings, not
# numeric cells). Columns D ("Price") and E ("Volume(1h)") contain values
# that look numeric/percentage, so Excel would normally coerce a plain
# .Value assignment into a real number (dropping significant trailing
# zeros, switching to scientific notation, etc.). To keep them as the
# literal text the diff expects, each such cell is forced to the "Text"
# number format ("@") immediately before its value is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

function Set-PlainCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Row 2 : BNB ---------------------------------------------------------
Set-TextCell "D2" "308.31"
Set-TextCell "E2" "-0.18%"

# --- Row 3 : OKB ----------------------------------------------------------
Set-TextCell "D3" "40.85"
Set-TextCell "E3" "1.94%"

# --- Row 4 : HuobiToken ----------------------------------------------------
Set-TextCell "D4" "5.132"
Set-TextCell "E4" "0.38%"

# --- Row 5 : Cronos --------------------------------------------------------
Set-TextCell "D5" "0.07626"
Set-TextCell "E5" "-1.36%"

# --- Row 6 : GateToken -> FTXToken -----------------------------------------
Set-PlainCell "B6" "FTXToken"
Set-PlainCell "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell  "D6" "1.605"
Set-TextCell  "E6" "-0.77%"

# --- Row 7 : FTXToken -> BTSEToken -----------------------------------------
Set-PlainCell "B7" "BTSEToken"
Set-PlainCell "C7" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell  "D7" "2.469"
Set-TextCell  "E7" "2.13%"

# --- Row 8 : BTSEToken -> MXToken ------------------------------------------
Set-PlainCell "B8" "MXToken"
Set-PlainCell "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell  "D8" "0.9043"
Set-TextCell  "E8" "2.55%"

# --- Row 9 : MXToken -> LiechtensteinCryptoassetsExchange ------------------
Set-PlainCell "B9" "LiechtensteinCryptoassetsExchange"
Set-PlainCell "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell  "D9" "0.1117"
Set-TextCell  "E9" "12.54%"

# --- Row 10 : LiechtensteinCryptoassetsExchange -> WazirX ------------------
Set-PlainCell "B10" "WazirX"
Set-PlainCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell  "D10" "0.1782"
Set-TextCell  "E10" "2.32%"

# --- Row 11 : WazirX -> MandalaExchangeToken -------------------------------
Set-PlainCell "B11" "MandalaExchangeToken"
Set-PlainCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell  "D11" "0.09186"
Set-TextCell  "E11" "1.74%"

# --- Row 12 : MandalaExchangeToken -> BitrueCoin ---------------------------
Set-PlainCell "B12" "BitrueCoin"
Set-PlainCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell  "D12" "0.04279"
Set-TextCell  "E12" "-3.69%"

# --- Row 13 : BitrueCoin -> BitMartToken -----------------------------------
Set-PlainCell "B13" "BitMartToken"
Set-PlainCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell  "D13" "0.1051"
Set-TextCell  "E13" "-0.31%"

# --- Row 14 : BitMartToken -> BitForexToken --------------------------------
Set-PlainCell "B14" "BitForexToken"
Set-PlainCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell  "D14" "0.001265"
Set-TextCell  "E14" "0.62%"

# --- Row 15 : BitForexToken -> TigerCash -----------------------------------
Set-PlainCell "B15" "TigerCash"
Set-PlainCell "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell  "D15" "0.005770"
Set-TextCell  "E15" "-1.43%"

# --- Row 16 : TigerCash -> LEO ----------------------------------------------
Set-PlainCell "B16" "LEO"
Set-PlainCell "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell  "D16" "3.349"
Set-TextCell  "E16" "-0.10%"

# --- Row 17 : LEO -> GateToken -----------------------------------------------
Set-PlainCell "B17" "GateToken"
Set-PlainCell "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell  "D17" "4.249"
Set-TextCell  "E17" "0.17%"

# --- Row 18 : BitpandaEcosystemToken (volume only) --------------------------
Set-TextCell "E18" "0.70%"

# --- Row 19 : MCDex -----------------------------------------------------------
Set-TextCell "D19" "6.648"
Set-TextCell "E19" "-5.47%"

# --- Row 20 : ProBitToken ------------------------------------------------------
Set-TextCell "D20" "0.1365"
Set-TextCell "E20" "1.12%"

# --- Row 21 : ZBToken ----------------------------------------------------------
Set-TextCell "D21" "0.2791"
Set-TextCell "E21" "-2.02%"

# --- Row 22 : CoinExToken -------------------------------------------------------
Set-TextCell "D22" "0.04080"
Set-TextCell "E22" "-0.90%"

# --- Row 23 : BitKan -------------------------------------------------------------
Set-TextCell "D23" "0.001238"
Set-TextCell "E23" "3.50%"

# --- Row 24 : HotbitToken ---------------------------------------------------------
Set-TextCell "D24" "0.004098"
Set-TextCell "E24" "-0.18%"

# --- Row 25 : NitroEx (volume only) -------------------------------------------------
Set-TextCell "E25" "-0.07%"

# --- Row 26 : UpBots (price only) ---------------------------------------------------
Set-TextCell "D26" "0.0003747"

# --- Row 38 : One ---------------------------------------------------------------------
Set-TextCell "D38" "0.02395"
Set-TextCell "E38" "1.70%"

# --- Row 39 : IDEX --------------------------------------------------------------------
Set-TextCell "D39" "0.05187"
Set-TextCell "E39" "-0.75%"

# --- Row 40 : KickToken ---------------------------------------------------------------
Set-TextCell "D40" "0.007781"
Set-TextCell "E40" "-2.18%"

# --- Row 41 : BKEXToken ---------------------------------------------------------------
Set-TextCell "D41" "0.1300"
Set-TextCell "E41" "-1.88%"

# --- Row 42 : Dexo --------------------------------------------------------------------
Set-TextCell "D42" "0.007053"
Set-TextCell "E42" "12.88%"

# --- Row 43 : CEJI (volume only) ------------------------------------------------------
Set-TextCell "E43" "-0.03%"

# --- Row 44 : LocalTraders ------------------------------------------------------------
Set-TextCell "D44" "0.007946"
Set-TextCell "E44" "-9.37%"

# --- Row 45 : PooCoin (volume only) ----------------------------------------------------
Set-TextCell "E45" "-7.62%"

# --- Row 46 : CoinLion -------------------------------------------------------------------
Set-TextCell "D46" "0.00006999"
Set-TextCell "E46" "6.49%"

# --- Row 47 : Kangarootoken ---------------------------------------------------------------
Set-TextCell "D47" "0.00000000751"
Set-TextCell "E47" "-0.07%"

# --- Row 48 : BOLO -------------------------------------------------------------------------
Set-TextCell "D48" "0.03154"
Set-TextCell "E48" "774.58%"

# --- Row 50 : CryptobidCoin -----------------------------------------------------------------
Set-TextCell "D50" "0.00002102"
Set-TextCell "E50" "-0.07%"

# --- Row 51 : SpecialPowerGold (volume only) -------------------------------------------------
Set-TextCell "E51" "-0.07%"
